$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "C.Tbilisi"

# 2. Row 6 (Urban) - some years become confidential/unavailable ("...")
$ws.Range("B6").Value = "..."
$ws.Range("C6").Value = "..."
$ws.Range("D6").Value = 156
$ws.Range("E6").Value = 128
$ws.Range("F6").Value = 143
$ws.Range("G6").Value = 143
$ws.Range("H6").Value = 159
$ws.Range("I6").Value = "..."
$ws.Range("J6").Value = 116
$ws.Range("K6").Value = 115
$ws.Range("L6").Value = 113
$ws.Range("M6").Value = 102
$ws.Range("N6").Value = 101
$ws.Range("O6").Value = 91

# 3. Row 7 (Rural) - some years become confidential/unavailable ("...")
$ws.Range("B7").Value = "..."
$ws.Range("C7").Value = "..."
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = "..."
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 11
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 5

# 4. New row 8 - footnote explaining the "..." placeholder, with a bold+underlined "Note:" lead-in
$noteCell = $ws.Range("A8")
$noteCell.Value = "Note: „ ... „ - Data is confidential or unavailable."
$noteCell.Font.Name = "Arial"
$noteCell.Font.Size = 9

$leadIn = $noteCell.Characters(1, 5)
$leadIn.Font.Name = "Arial"
$leadIn.Font.Size = 9
$leadIn.Font.Bold = $true
$leadIn.Font.Underline = $true

$rest = $noteCell.Characters(6, 47)
$rest.Font.Name = "Arial"
$rest.Font.Size = 9
